$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for row 2: category, title, and refugee weight columns
$ws.Range("D2").Value = "Refugee"
$ws.Range("E2").Value = "Refugee Travel Documents in Poland"
$ws.Range("M2").Value = 20

# Match the author's final selection/view state (row 2 scrolled into view, K2 active)
$ws.Range("K2").Select()
$excel.ActiveWindow.ScrollRow = 2
